# Apply the Alvearie FHIR IG "source-data-model-version" StructureDefinition
# regeneration edits (Version bump 5.0.0 -> 6.0.0, Date refresh, Publisher and
# Jurisdiction details added in place of the placeholder "Contact" rows, and
# the generated root Extension's Short/Definition text updated) to the
# Metadata and Elements worksheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Metadata" worksheet (sheet 1): key/value property table
# ---------------------------------------------------------------------------
$meta = $wb.Worksheets.Item(1)

# The source sheet has two identical placeholder "Contact" rows (10 and 11).
# Remove the second one; everything below shifts up by one row, which is why
# the sheet dimension goes from A1:B21 to A1:B20.
$meta.Range("A11").EntireRow.Delete()

# Version: 5.0.0 -> 6.0.0
$meta.Range("B3").Value = "6.0.0"

# Date refreshed to the new publication timestamp
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value filled in
$meta.Range("B9").Value = "Alvearie Team"

# The remaining (former second) "Contact" row is repurposed for Jurisdiction
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ---------------------------------------------------------------------------
# "Elements" worksheet (sheet 2): StructureDefinition element table
# ---------------------------------------------------------------------------
$elements = $wb.Worksheets.Item(2)

# Root Extension row (row 2): Short / Definition now reflect the real
# title/description instead of the generic "Extension" / "An Extension".
$elements.Range("K2").Value = "Source Data Model Version"
$elements.Range("L2").Value = "Version of the source system's data model, used by either the data producer or the data integrator"
